# Updated cryptos list on Fri Jul  7 10:35:00 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row on Sheet1 with the latest scrape. Both columns hold plain
# text in the source data (prices use "." as a thousands separator in
# several rows, e.g. "30.123.74", and the volume values keep their
# surrounding padding spaces, e.g. "  -3.29%  "), so this must land back in
# the workbook as literal text -- not get auto-coerced into a number by
# Excel's input parser.
#
# To guarantee that, each "Price" cell is briefly switched to the Text
# number format ("@") before the new value is written, then the format is
# reset to General and the cell's Style is reset to "Normal" so the
# produced workbook keeps using the original default cell style (no new
# styles get left behind on the cell itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.123.74'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.29'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.61'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.32%  '

$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4650'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2817'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06538'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.88'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07813'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.11'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.858.23'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.115'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6670'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '280.26'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.151.58'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.21%  '

$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.497'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.59'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.101.66'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007234'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.130'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.306'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.34'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.82'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.98%  '

$ws.Range("E28").Value = '  -9.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.342'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09578'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.393'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.467'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.099'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04644'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6978'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.090'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.34%  '

$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01850'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.279'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.32%  '

$ws.Range("E40").Value = '  -4.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.87'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8518'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.916'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.45%  '

$ws.Range("E44").Value = '  -0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4146'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.24'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '989.82'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.160'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.222'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.03'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1137'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.99%  '
